$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70: B70 was stored as a text "4" (inline string); change it to a real number 4
$ws.Range("B70").Value = 4

# Append a new row 71 with the new annotation data
$ws.Range("A71").Value = "Ruilin"

# B71 keeps the numeric-looking value "3" stored as text (matches source data)
$ws.Range("B71").NumberFormat = "@"
$ws.Range("B71").Value = "3"
$ws.Range("B71").Style = "Normal"

$ws.Range("C71").Value = "无"
$ws.Range("D71").Value = "DIS"
$ws.Range("E71").Value = "OTH"
$ws.Range("F71").Value = "3c70bee3-3ebe-492b-b68a-cb43e1a99f35"
$ws.Range("G71").Value = "H1Ww66x0-_annotated.xlsx"
$ws.Range("H71").Value = "We will include additional details on the hyper-parameters of the baselines for clarity."
